$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.471.53"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "1.919.56"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("D4").Value = "'1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.55%  "
$ws.Range("D5").Value = "'325.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").Value = "'1.008"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.63%  "
$ws.Range("D7").Value = "'0.4828"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.11%  "
$ws.Range("D8").Value = "'0.4080"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.53%  "
$ws.Range("D9").Value = "'0.08236"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.32%  "
$ws.Range("D10").Value = "'1.017"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.56%  "
$ws.Range("D11").Value = "'23.46"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.72%  "
$ws.Range("D12").Value = "1.906.76"
$ws.Range("E12").Value = "  +1.51%  "
$ws.Range("D13").Value = "'6.086"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.44%  "
$ws.Range("D14").Value = "'7.254"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.71%  "
$ws.Range("E15").Value = "  +2.02%  "
$ws.Range("D16").Value = "'0.06809"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.13%  "
$ws.Range("D17").Value = "'1.008"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.69%  "
$ws.Range("E18").Value = "  +1.14%  "
$ws.Range("D19").Value = "'17.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("E20").Value = "  +0.58%  "
$ws.Range("D21").Value = "29.495.43"
$ws.Range("E21").Value = "  +0.40%  "
$ws.Range("D22").Value = "'5.662"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.52%  "
$ws.Range("D23").Value = "'11.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.22%  "
$ws.Range("D24").Value = "'2.180"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.08%  "
$ws.Range("D25").Value = "2.141.95"
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("D26").Value = "'6.663"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +10.44%  "
$ws.Range("D27").Value = "'156.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.01%  "
$ws.Range("E28").Value = "  +1.69%  "
$ws.Range("D29").Value = "'2.118"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.50%  "
$ws.Range("E30").Value = "  +2.08%  "
$ws.Range("D31").Value = "'1.022"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("D32").Value = "'0.09605"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.19%  "
$ws.Range("D33").Value = "'5.682"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.82%  "
$ws.Range("E34").Value = "  +0.74%  "
$ws.Range("D35").Value = "'1.374"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.76%  "
$ws.Range("D36").Value = "'0.02290"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.92%  "
$ws.Range("D37").Value = "'0.06114"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.13%  "
$ws.Range("E38").Value = "  +0.95%  "
$ws.Range("D39").Value = "'8.098"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.95%  "
$ws.Range("D40").Value = "'0.5998"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.38%  "
$ws.Range("D41").Value = "'10.88"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.44%  "
$ws.Range("D42").Value = "'0.1851"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.61%  "
$ws.Range("D43").Value = "'2.412"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("D44").Value = "'1.280"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.62%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'12.48"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.64%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").Value = "'0.07601"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.49%  "
$ws.Range("D47").Value = "'0.5594"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.74%  "
$ws.Range("D48").Value = "'1.958"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.08%  "
$ws.Range("D49").Value = "'118.56"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'2.427"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.06%  "
$ws.Range("D51").Value = "'72.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.98%  "
